$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.228.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.05%  "
$ws.Range("D3").Value = "'1.814.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.43%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'329.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.4416"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("D8").Value = "'0.3700"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.11%  "
$ws.Range("D9").Value = "'44.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "'0.07708"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.48%  "
$ws.Range("D11").Value = "'1.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'22.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Value = "'7.571"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.53%  "
$ws.Range("D15").Value = "'6.243"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "'1.823.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.13%  "
$ws.Range("D17").Value = "'92.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.31%  "
$ws.Range("D18").Value = "'0.00001084"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.34%  "
$ws.Range("D19").Value = "'0.06569"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.67%  "
$ws.Range("D20").Value = "'0.9994"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "'17.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.29%  "
$ws.Range("D22").Value = "'6.209"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.13%  "
$ws.Range("D23").Value = "'28.296.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("D24").Value = "'11.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.83%  "
$ws.Range("D25").Value = "'1.996"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -16.71%  "
$ws.Range("D26").Value = "'20.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.68%  "
$ws.Range("D27").Value = "'156.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.18%  "
$ws.Range("D28").Value = "'2.015.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.20%  "
$ws.Range("D29").Value = "'2.318"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").Value = "'127.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.37%  "
$ws.Range("D31").Value = "'1.202"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.21%  "
$ws.Range("D32").Value = "'5.872"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.15%  "
$ws.Range("D33").Value = "'0.09225"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.62%  "
$ws.Range("D34").Value = "'3.667"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").Value = "'13.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.88%  "
$ws.Range("D36").Value = "'0.02353"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").Value = "'0.2169"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").Value = "'5.169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("D39").Value = "'0.6583"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("D40").Value = "'0.06197"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").Value = "'1.196"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "'8.108"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.55%  "
$ws.Range("D43").Value = "'0.9991"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").Value = "'13.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.38%  "
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").Value = "'0.6076"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.33%  "
$ws.Range("D47").Value = "'3.760"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").Value = "'126.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").Value = "'2.036"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.21%  "
$ws.Range("D50").Value = "'1.153"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.55%  "
$ws.Range("D51").Value = "'0.06979"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.43%  "

Write-Host "DONE"